$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing A:T data to B:U
$ws.Range("A1").EntireColumn.Insert()

# Set the new column's header and narrower width
$ws.Range("A1").Value2 = "Lab. #"
$ws.Columns("A:A").ColumnWidth = 6.75

# Fill in the laboratory numbers for each data row
$ws.Range("A2").Value2 = 10815
$ws.Range("A3").Value2 = 10989
$ws.Range("A4").Value2 = 10815
$ws.Range("A5").Value2 = 10990
$ws.Range("A6").Value2 = 10815
$ws.Range("A7").Value2 = 10991
$ws.Range("A8").Value2 = 10815
$ws.Range("A9").Value2 = 10992
$ws.Range("A10").Value2 = 10815
$ws.Range("A11").Value2 = 10993
$ws.Range("A12").Value2 = 10815

# Highlight every other data row (lab number 10815 repeat group) with a light green fill
$ws.Range("A2:U2").Interior.Color = 12379352
$ws.Range("A4:U4").Interior.Color = 12379352
$ws.Range("A6:U6").Interior.Color = 12379352
$ws.Range("A8:U8").Interior.Color = 12379352
$ws.Range("A10:U10").Interior.Color = 12379352
$ws.Range("A12:U12").Interior.Color = 12379352
